$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-03-07 20:35:20"

for ($row = 2; $row -le 51; $row++) {
    $ws.Cells.Item($row, 8).Value = $newTimestamp
}
